$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $oldText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $ptext = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($ptext -eq $oldText) {
            return $i
        }
    }
    return -1
}

# Version/date line
$idx = Find-ParagraphIndex $d 'Version Date: Tue Aug 28 12:33:50 2018 -0700'
if ($idx -eq -1) { throw "paragraph not found: Version/date line" }
$d.Paragraphs.Item($idx).Range.Text = 'Version Date: Wed Aug 29 17:32:15 2018 -0700'

# Introduction paragraph 1
$idx = Find-ParagraphIndex $d 'Bioassessment is an essential element of environmental monitoring programs that informs decisions for managing aquatic resources. Decades of research have supported the development of methods that use a variety of assemblages with regional applications in streams, rivers, lakes, and marine environments. This body of applied tools represents significant achievements in developing accurate and interpretable assessment methods that rely on biological organisms as sentinels of environmental condition. Monitoring programs in the United States and internationally have collected millions of records of biological data spanning decades and hundreds of assessment methods have been developed from these data, yet the ability of these tools to reach the management community and to positively affect environmental change is imbalanced relative to the abundance of information that could be used. Morever, existing methods can lack transparency, require specialized training to implement, and are often not discoverable beyond specific research applications. Decision-makers require additional tools that synthesize information and bridge the gap between method and application.'
if ($idx -eq -1) { throw "paragraph not found: Introduction paragraph 1" }
$d.Paragraphs.Item($idx).Range.Text = 'Bioassessment is an essential element of environmental monitoring programs that informs decisions for managing aquatic resources. Decades of research have supported the development of methods that use a variety of assemblages with regional applications in streams, rivers, lakes, and marine environments. This body of applied tools represents significant achievements in developing accurate and interpretable assessment methods that rely on biological organisms as sentinels of environmental condition. Monitoring programs in the United States and internationally have collected millions of records of biological data spanning decades and hundreds of assessment methods have been developed from these data, yet the ability of these tools to reach the management community and to positively affect environmental change is imbalanced relative to the amount of information that is available through research and coordinated monitoring efforts. Morever, existing methods can lack transparency, require specialized training to implement, and are often not discoverable beyond specific research applications. Decision-makers require additional tools that synthesize information and bridge the gap between method and application.'

# Introduction paragraph 2
$idx = Find-ParagraphIndex $d 'Legal mandates to assess biological condition have set a precedent for developing bioassessment methods in the United States, Canada, and Europe. Since the passage of these mandates, the science of bioassessment has for decades focused on addressing technical challenges for developing robust indices that accurately describe environmental condition. Basic research questions to address these goals focused on identifying which biological components respond predictably to environmental change, how these components can be measured with minimum uncerainty, and what basis of comparison is used to evaluate relative changes between communities. Most bioassessment indices developed at the assemblage-level are characterized as either multimetric, such as the index of biotic integrity, or multivariate where condition is based on ratios between observed and expected taxa. The reference-condition approach also establishes the foundation for many bioassessment methods whereby a set of reference sites are identified and used to evaluate accepted levels of deviation in scores to define potential impacts.'
if ($idx -eq -1) { throw "paragraph not found: Introduction paragraph 2" }
$d.Paragraphs.Item($idx).Range.Text = 'The science of bioassessment has for decades focused on addressing technical challenges for developing indices that accurately describe environmental condition. Legal mandates to assess biological condition have set a precedent for developing bioassessment methods in the United States (Clean Water Act), Canada (Canada Waters Act), and Europe (Water Framework Directive). Basic research to address broad legislative objectives has historically focused on identifying which biological components respond predictably to environmental change, how these components can be measured with minimum uncerainty, and what basis of comparison is used to evaluate relative changes between communities. Most bioassessment indices developed at the assemblage-level are characterized as either multimetric, such as the index of biotic integrity, or multivariate where condition is assessed using ratios of observed and expected taxa. The reference-condition approach also establishes the foundation for many bioassessment methods whereby a set of reference sites are identified and used to evaluate levels of biological deviation to define potential impacts.'

# Introduction paragraph 3
$idx = Find-ParagraphIndex $d 'Methods for evaluating and comparing biological communities that have been developed by the research community have generally been accepted as reliable means for developing bioassessment tools. As a result, most assessment tools follow a similar approach to development that is often guided by technical support documents that summarize the body of research to date on bioassessment. This has in part contributed to the proliferation of hundreds of assessment methods that have been developed for specic regional applcations. Although there are logistical rationale for why location-specific methods are needed, some have argued that larger-scale and more integrated methods are more practical given the obvious benefit of assessment standards that have broader applicability (Whitter study). What is the point of this paragraph?'
if ($idx -eq -1) { throw "paragraph not found: Introduction paragraph 3" }
$d.Paragraphs.Item($idx).Range.Text = 'Methods for evaluating and comparing biological communities that have been developed by the research community have generally been accepted by environmental managers as robust for developing bioassessment tools. As a result, new assessment tools can be developed using readily available technical support documents that summarize the body of research and best practices to date on bioassessment. This has in part contributed to the proliferation of hundreds of assessment methods that have been developed for specific regional applcations. Although there are logistical and ecological rationale for why location- and taxa-specific methods are needed, concerns about redundancy, duplicated effort, and lack of coordinated monitoring have recently been discussed within the research community. Morever, the abundance of available methods can be a point of frustration for managers given a lack of guidance for choosing an appropriate method among alternatives. The process to charactize how an index could be used in practice to inform decisions and prioritize management actions is often opaque relative to why an index may have originally been developed.'

# Introduction paragraph 4
$idx = Find-ParagraphIndex $d 'What distinguishes bioassessment from others field of study in ecology is the link to environmental management. Although bioassessment methods can and have been used to inform basic research, the intended use of these tools, often by legal definition, is in the applied sense. An argument could be made that an index is only as valuable as its level of integration with the management and regulatory communities. In the United States, the CWA gives the power to states, tribes, and territories to delop their own methods which in turn require approval at the federal level to be implemented into a regulatory framework, e.g., TMDL reporting, permitting, etc.. If federal approval is a rough assessment of index efficacy, a tremendous imbalance exists between the methods developed and those are officially used in a regulatory context.'
if ($idx -eq -1) { throw "paragraph not found: Introduction paragraph 4" }
$d.Paragraphs.Item($idx).Range.Text = 'What distinguishes bioassessment from others field of ecological study is the link to environmental management. Although bioassessment methods can and have been used to inform basic research, the intended use of these tools, often through legal definition, is in the applied sense. An argument could be made that an index is only as valuable as its level of integration with management and regulatory communities. In the United States, the CWA gives the power to states, tribes, and territories to develop their own methods, which in turn require federal approval to be implemented into a regulatory framework, e.g., TMDL reporting, permitting, etc. If federal approval is a rough assessment of index efficacy, a tremendous imbalance exists between the methods developed and those that are federally approved for regulatory use. Of those that are approved, even more concerning is the manner of application within standard regulatory frameworks. Biological indices are typically used to develop post-hoc diagnoses that can trigger remediative or restoration actions. A growing concern is that these tools, although technically sound, are primarliy being used to document the long-term demise of environmental health. A much broader use for bioassessment to pro-actively guide planning decisions, such as identifying conservation priorities, could greatly extend the reach of tools that have already been developed.'

# Final paragraph: replace its text and split off a new trailing paragraph (same BodyText style)
# carrying the second block of new text.
$idx = Find-ParagraphIndex $d 'What is the reason for this imbalance?'
if ($idx -eq -1) { throw "paragraph not found: final imbalance paragraph" }

# Setting text with a trailing carriage return both replaces the original paragraph
# text and splits off a new (empty) paragraph immediately after it.
$d.Paragraphs.Item($idx).Range.Text = 'Bioassessment currently suffers from an excess of information and forward progress will not be made unless this information meets the needs of the management community. A new mode of operation is needed whereby method development is open and transparent, existing methods are discoverable and reproducible, and information transfer to the management community is intuitive and purposeful. Open science principles that can democratize all aspects of the scientific method can meet these needs, yet bioassessment research and its application to better serve the environment has not fully embraced these principles. Others have advocated more broadly for inclusion of these principles in the ecological sciences (Hampton et al. 2015) and a growing wave of momentum has seen open science permeate how scientists conceptualize research in other disciplines (e.g., archeaology, behavioral ecology, vegetation sciences). Adopting an open science paradigm in biaossessment is particularly relevant compared to other fields given the explicit need to develop tools that are open and accessible to the management community. Legal and even ethical precedents in bioassessment may necessitate the open sharing of data given that environmental monitoring programs are often publicly funded.' + "`r"

# The freshly split-off paragraph starts with no runs of its own; writing a trailing
# space first (then immediately overwriting without it) matches the whitespace-run
# formatting used by every other paragraph in this document.
$d.Paragraphs.Item($idx + 1).Range.Text = 'This review will empower the research and management community to embrace open science as a new mode of thinking for bioassessment. These approaches are expected to benefit the research community by augmenting existing workflows for developing assessment tools, but more importantly, improve the ability of these methods to address environmental issues by bridging the gap between the scientific, management, and regulatory communities. An overview of the general principles of open science is provided, followed by a discussion of specific benefits and how these principles can be applied to bioassessment. We use examples from the state of California to demonstrate a real world scenario of how existing tools can be tailored to address legislative mandates for free and open sharing of data. We conclude with a discussion of technical, sociocultural, and instutional hurdles that have thus far prevented widespread adoption of opens science and provide recommendations for the bioassessment community to address these challenges.' + " "
$d.Paragraphs.Item($idx + 1).Range.Text = 'This review will empower the research and management community to embrace open science as a new mode of thinking for bioassessment. These approaches are expected to benefit the research community by augmenting existing workflows for developing assessment tools, but more importantly, improve the ability of these methods to address environmental issues by bridging the gap between the scientific, management, and regulatory communities. An overview of the general principles of open science is provided, followed by a discussion of specific benefits and how these principles can be applied to bioassessment. We use examples from the state of California to demonstrate a real world scenario of how existing tools can be tailored to address legislative mandates for free and open sharing of data. We conclude with a discussion of technical, sociocultural, and instutional hurdles that have thus far prevented widespread adoption of opens science and provide recommendations for the bioassessment community to address these challenges.'

Write-Output "edits applied"
